$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 69.89967366666667
$ws.Cells.Item(2, 8).Value = 209.699021
$ws.Cells.Item(2, 9).Value = 0.6608367681537789
$ws.Cells.Item(2, 10).Value = 0.660836768153779
$ws.Cells.Item(2, 13).Value = 4.717738333333333
$ws.Cells.Item(2, 14).Value = 14.153215
$ws.Cells.Item(2, 15).Value = 0.2002263444295212
$ws.Cells.Item(2, 16).Value = 0.2002263444295212
$ws.Cells.Item(2, 17).Value = 329.7683699447239
$ws.Cells.Item(2, 18).Value = 2967.915329502515
$ws.Cells.Item(2, 19).Value = 0.1323169303520502
$ws.Cells.Item(2, 20).Value = 0.1323169303520502

$ws.Cells.Item(3, 7).Value = 69.89967366666667
$ws.Cells.Item(3, 8).Value = 209.699021
$ws.Cells.Item(3, 9).Value = 0.6608367681537789
$ws.Cells.Item(3, 10).Value = 0.660836768153779
$ws.Cells.Item(3, 14).Value = 9.228847
$ws.Cells.Item(3, 15).Value = 0.1305610278731266
$ws.Cells.Item(3, 16).Value = 0.1305610278731266
$ws.Cells.Item(3, 17).Value = 215.0311312065319
$ws.Cells.Item(3, 18).Value = 1935.280180858787
$ws.Cells.Item(3, 19).Value = 0.08627952770651243
$ws.Cells.Item(3, 20).Value = 0.08627952770651244

$ws.Cells.Item(4, 7).Value = 69.89967366666667
$ws.Cells.Item(4, 8).Value = 209.699021
$ws.Cells.Item(4, 9).Value = 0.6608367681537789
$ws.Cells.Item(4, 10).Value = 0.660836768153779
$ws.Cells.Item(4, 13).Value = 0.6908423333333333
$ws.Cells.Item(4, 14).Value = 2.072527
$ws.Cells.Item(4, 15).Value = 0.02932015834857891
$ws.Cells.Item(4, 16).Value = 0.02932015834857891
$ws.Cells.Item(4, 17).Value = 48.28965365511856
$ws.Cells.Item(4, 18).Value = 434.606882896067
$ws.Cells.Item(4, 19).Value = 0.01937583868483193
$ws.Cells.Item(4, 20).Value = 0.01937583868483193

$ws.Cells.Item(5, 7).Value = 69.89967366666667
$ws.Cells.Item(5, 8).Value = 209.699021
$ws.Cells.Item(5, 9).Value = 0.6608367681537789
$ws.Cells.Item(5, 10).Value = 0.660836768153779
$ws.Cells.Item(5, 13).Value = 15.077163
$ws.Cells.Item(5, 14).Value = 45.231489
$ws.Cells.Item(5, 15).Value = 0.6398924693487733
$ws.Cells.Item(5, 16).Value = 0.6398924693487733
$ws.Cells.Item(5, 17).Value = 1053.888773519141
$ws.Cells.Item(5, 18).Value = 9484.998961672269
$ws.Cells.Item(5, 19).Value = 0.4228644714103844
$ws.Cells.Item(5, 20).Value = 0.4228644714103845

$ws.Cells.Item(6, 9).Value = 0.1661491941864736
$ws.Cells.Item(6, 10).Value = 0.1661491941864736
$ws.Cells.Item(6, 13).Value = 4.717738333333333
$ws.Cells.Item(6, 14).Value = 14.153215
$ws.Cells.Item(6, 15).Value = 0.2002263444295212
$ws.Cells.Item(6, 16).Value = 0.2002263444295212
$ws.Cells.Item(6, 17).Value = 82.91116895262222
$ws.Cells.Item(6, 18).Value = 746.2005205736
$ws.Cells.Item(6, 19).Value = 0.03326744578186826
$ws.Cells.Item(6, 20).Value = 0.03326744578186827

$ws.Cells.Item(7, 9).Value = 0.1661491941864736
$ws.Cells.Item(7, 10).Value = 0.1661491941864736
$ws.Cells.Item(7, 14).Value = 9.228847
$ws.Cells.Item(7, 15).Value = 0.1305610278731266
$ws.Cells.Item(7, 16).Value = 0.1305610278731266
$ws.Cells.Item(7, 17).Value = 54.06365217054223
$ws.Cells.Item(7, 18).Value = 486.57286953488
$ws.Cells.Item(7, 19).Value = 0.0216926095732777
$ws.Cells.Item(7, 20).Value = 0.02169260957327771

$ws.Cells.Item(8, 9).Value = 0.1661491941864736
$ws.Cells.Item(8, 10).Value = 0.1661491941864736
$ws.Cells.Item(8, 13).Value = 0.6908423333333333
$ws.Cells.Item(8, 14).Value = 2.072527
$ws.Cells.Item(8, 15).Value = 0.02932015834857891
$ws.Cells.Item(8, 16).Value = 0.02932015834857891
$ws.Cells.Item(8, 17).Value = 12.14110265800889
$ws.Cells.Item(8, 18).Value = 109.26992392208
$ws.Cells.Item(8, 19).Value = 0.004871520683036193
$ws.Cells.Item(8, 20).Value = 0.004871520683036194

$ws.Cells.Item(9, 9).Value = 0.1661491941864736
$ws.Cells.Item(9, 10).Value = 0.1661491941864736
$ws.Cells.Item(9, 13).Value = 15.077163
$ws.Cells.Item(9, 14).Value = 45.231489
$ws.Cells.Item(9, 15).Value = 0.6398924693487733
$ws.Cells.Item(9, 16).Value = 0.6398924693487733
$ws.Cells.Item(9, 17).Value = 264.97128931184
$ws.Cells.Item(9, 18).Value = 2384.74160380656
$ws.Cells.Item(9, 19).Value = 0.1063176181482915
$ws.Cells.Item(9, 20).Value = 0.1063176181482915

$ws.Cells.Item(10, 7).Value = 4.152730666666667
$ws.Cells.Item(10, 8).Value = 12.458192
$ws.Cells.Item(10, 9).Value = 0.0392602278210887
$ws.Cells.Item(10, 10).Value = 0.03926022782108871
$ws.Cells.Item(10, 13).Value = 4.717738333333333
$ws.Cells.Item(10, 14).Value = 14.153215
$ws.Cells.Item(10, 15).Value = 0.2002263444295212
$ws.Cells.Item(10, 16).Value = 0.2002263444295212
$ws.Cells.Item(10, 17).Value = 19.59149665414222
$ws.Cells.Item(10, 18).Value = 176.32346988728
$ws.Cells.Item(10, 19).Value = 0.007860931898086775
$ws.Cells.Item(10, 20).Value = 0.007860931898086778

$ws.Cells.Item(11, 7).Value = 4.152730666666667
$ws.Cells.Item(11, 8).Value = 12.458192
$ws.Cells.Item(11, 9).Value = 0.0392602278210887
$ws.Cells.Item(11, 10).Value = 0.03926022782108871
$ws.Cells.Item(11, 14).Value = 9.228847
$ws.Cells.Item(11, 15).Value = 0.1305610278731266
$ws.Cells.Item(11, 16).Value = 0.1305610278731266
$ws.Cells.Item(11, 17).Value = 12.77497198495822
$ws.Cells.Item(11, 18).Value = 114.974747864624
$ws.Cells.Item(11, 19).Value = 0.005125855698854462
$ws.Cells.Item(11, 20).Value = 0.005125855698854462

$ws.Cells.Item(12, 7).Value = 4.152730666666667
$ws.Cells.Item(12, 8).Value = 12.458192
$ws.Cells.Item(12, 9).Value = 0.0392602278210887
$ws.Cells.Item(12, 10).Value = 0.03926022782108871
$ws.Cells.Item(12, 13).Value = 0.6908423333333333
$ws.Cells.Item(12, 14).Value = 2.072527
$ws.Cells.Item(12, 15).Value = 0.02932015834857891
$ws.Cells.Item(12, 16).Value = 0.02932015834857891
$ws.Cells.Item(12, 17).Value = 2.868882143464889
$ws.Cells.Item(12, 18).Value = 25.819939291184
$ws.Cells.Item(12, 19).Value = 0.001151116096515604
$ws.Cells.Item(12, 20).Value = 0.001151116096515604

$ws.Cells.Item(13, 7).Value = 4.152730666666667
$ws.Cells.Item(13, 8).Value = 12.458192
$ws.Cells.Item(13, 9).Value = 0.0392602278210887
$ws.Cells.Item(13, 10).Value = 0.03926022782108871
$ws.Cells.Item(13, 13).Value = 15.077163
$ws.Cells.Item(13, 14).Value = 45.231489
$ws.Cells.Item(13, 15).Value = 0.6398924693487733
$ws.Cells.Item(13, 16).Value = 0.6398924693487733
$ws.Cells.Item(13, 17).Value = 62.611397156432
$ws.Cells.Item(13, 18).Value = 563.502574407888
$ws.Cells.Item(13, 19).Value = 0.02512232412763186
$ws.Cells.Item(13, 20).Value = 0.02512232412763186

$ws.Cells.Item(14, 7).Value = 14.14774133333333
$ws.Cells.Item(14, 8).Value = 42.443224
$ws.Cells.Item(14, 9).Value = 0.1337538098386587
$ws.Cells.Item(14, 10).Value = 0.1337538098386588
$ws.Cells.Item(14, 13).Value = 4.717738333333333
$ws.Cells.Item(14, 14).Value = 14.153215
$ws.Cells.Item(14, 15).Value = 0.2002263444295212
$ws.Cells.Item(14, 16).Value = 0.2002263444295212
$ws.Cells.Item(14, 17).Value = 66.74534161835112
$ws.Cells.Item(14, 18).Value = 600.70807456516
$ws.Cells.Item(14, 19).Value = 0.02678103639751596
$ws.Cells.Item(14, 20).Value = 0.02678103639751597

$ws.Cells.Item(15, 7).Value = 14.14774133333333
$ws.Cells.Item(15, 8).Value = 42.443224
$ws.Cells.Item(15, 9).Value = 0.1337538098386587
$ws.Cells.Item(15, 10).Value = 0.1337538098386588
$ws.Cells.Item(15, 14).Value = 9.228847
$ws.Cells.Item(15, 15).Value = 0.1305610278731266
$ws.Cells.Item(15, 16).Value = 0.1305610278731266
$ws.Cells.Item(15, 17).Value = 43.52244672030312
$ws.Cells.Item(15, 18).Value = 391.702020482728
$ws.Cells.Item(15, 19).Value = 0.017463034894482
$ws.Cells.Item(15, 20).Value = 0.017463034894482

$ws.Cells.Item(16, 7).Value = 14.14774133333333
$ws.Cells.Item(16, 8).Value = 42.443224
$ws.Cells.Item(16, 9).Value = 0.1337538098386587
$ws.Cells.Item(16, 10).Value = 0.1337538098386588
$ws.Cells.Item(16, 13).Value = 0.6908423333333333
$ws.Cells.Item(16, 14).Value = 2.072527
$ws.Cells.Item(16, 15).Value = 0.02932015834857891
$ws.Cells.Item(16, 16).Value = 0.02932015834857891
$ws.Cells.Item(16, 17).Value = 9.773858634116445
$ws.Cells.Item(16, 18).Value = 87.964727707048
$ws.Cells.Item(16, 19).Value = 0.003921682884195186
$ws.Cells.Item(16, 20).Value = 0.003921682884195187

$ws.Cells.Item(17, 7).Value = 14.14774133333333
$ws.Cells.Item(17, 8).Value = 42.443224
$ws.Cells.Item(17, 9).Value = 0.1337538098386587
$ws.Cells.Item(17, 10).Value = 0.1337538098386588
$ws.Cells.Item(17, 13).Value = 15.077163
$ws.Cells.Item(17, 14).Value = 45.231489
$ws.Cells.Item(17, 15).Value = 0.6398924693487733
$ws.Cells.Item(17, 16).Value = 0.6398924693487733
$ws.Cells.Item(17, 17).Value = 213.307802164504
$ws.Cells.Item(17, 18).Value = 1919.770219480536
$ws.Cells.Item(17, 19).Value = 0.08558805566246558
$ws.Cells.Item(17, 20).Value = 0.08558805566246561
